$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = [double]"1.0000000004200407E-4"
$ws.Range("B1").Value = [double]"2.5999199066825449E-2"
$ws.Range("C1").Value = [double]"4.5337094612479381E-2"
$ws.Range("D1").Value = [double]"150051042.74125579"
$ws.Range("E1").Value = [double]"30.200738390273099"
$ws.Range("F1").Value = [double]"52.429561697718796"
$ws.Range("G1").Value = [double]"2.3503953926105909E-14"
$ws.Range("H1").Value = [double]"0.59544112184144316"
$ws.Range("I1").Value = [double]"6.1535020581943885E-3"
$ws.Range("J1").Value = [double]"6.772723344689259E-6"
$ws.Range("K1").Value = [double]"0.96949840104266327"
$ws.Range("L1").Value = [double]"-0.49306638370346012"
$ws.Range("M1").Value = [double]"-0.69928252362742449"
$ws.Range("A2").Value = [double]"7.2349937893170079E-3"
$ws.Range("B2").Value = [double]"0.31016695466413624"
$ws.Range("C2").Value = [double]"4.8616542313354626E-2"
$ws.Range("D2").Value = [double]"254800367.08372283"
$ws.Range("E2").Value = [double]"30.025584739404472"
$ws.Range("F2").Value = [double]"60.255254748468253"
$ws.Range("G2").Value = [double]"0.8753582519130122"
$ws.Range("H2").Value = [double]"0.26866431961982734"
$ws.Range("I2").Value = [double]"4.2026759163526443E-6"
$ws.Range("J2").Value = [double]"3.483585410102172E-5"
$ws.Range("K2").Value = [double]"0.98257186998833945"
$ws.Range("L2").Value = [double]"0.92755812393150383"
$ws.Range("M2").Value = [double]"0.51595613224007097"
$ws.Range("A3").Value = [double]"1.0000000002225827E-4"
$ws.Range("B3").Value = [double]"1.0165936996063066E-2"
$ws.Range("C3").Value = [double]"4.5445067532555054E-2"
$ws.Range("D3").Value = [double]"150000105.93788049"
$ws.Range("E3").Value = [double]"36.272455679572722"
$ws.Range("F3").Value = [double]"52.287946926810186"
$ws.Range("G3").Value = [double]"2.3330539731932203E-14"
$ws.Range("H3").Value = [double]"0.59999685189836072"
$ws.Range("I3").Value = [double]"8.1776069520533361E-6"
$ws.Range("J3").Value = [double]"2.1550843287448077E-3"
$ws.Range("K3").Value = [double]"0.77076967697899634"
$ws.Range("L3").Value = [double]"-0.53911458037838278"
$ws.Range("M3").Value = [double]"-0.73657496101126596"
$ws.Range("A4").Value = [double]"1.0000000002231741E-4"
$ws.Range("B4").Value = [double]"1.0208397609236719E-2"
$ws.Range("C4").Value = [double]"4.5433941259353523E-2"
$ws.Range("D4").Value = [double]"150000310.13718903"
$ws.Range("E4").Value = [double]"42.080803172244494"
$ws.Range("F4").Value = [double]"57.063889827745115"
$ws.Range("G4").Value = [double]"2.4864129457870661E-14"
$ws.Range("H4").Value = [double]"0.59999032808394559"
$ws.Range("I4").Value = [double]"1.5002641406479561E-5"
$ws.Range("J4").Value = [double]"2.2040099058255732E-3"
$ws.Range("K4").Value = [double]"0.94096482201741871"
$ws.Range("L4").Value = [double]"-0.51308744140307105"
$ws.Range("M4").Value = [double]"-0.55867789773326226"
$ws.Range("A5").Value = [double]"1.0000000002224167E-4"
$ws.Range("B5").Value = [double]"1.0358933111103337E-2"
$ws.Range("C5").Value = [double]"4.5455777291659262E-2"
$ws.Range("D5").Value = [double]"150000064.36837882"
$ws.Range("E5").Value = [double]"35.429475836197291"
$ws.Range("F5").Value = [double]"51.282208357383311"
$ws.Range("G5").Value = [double]"2.2959368026583736E-14"
$ws.Range("H5").Value = [double]"0.59999708593979983"
$ws.Range("I5").Value = [double]"8.7887515596040539E-7"
$ws.Range("J5").Value = [double]"2.1427824619856097E-3"
$ws.Range("K5").Value = [double]"0.66943611717653306"
$ws.Range("L5").Value = [double]"-0.50931480221125791"
$ws.Range("M5").Value = [double]"-0.57018565567069657"
$ws.Range("A6").Value = [double]"1.0000000002233415E-4"
$ws.Range("B6").Value = [double]"1.0354160372415875E-2"
$ws.Range("C6").Value = [double]"4.5423751631699677E-2"
$ws.Range("D6").Value = [double]"150000130.72154206"
$ws.Range("E6").Value = [double]"37.223930590336622"
$ws.Range("F6").Value = [double]"51.653729351762557"
$ws.Range("G6").Value = [double]"2.5370065585859287E-14"
$ws.Range("H6").Value = [double]"0.59998906543721109"
$ws.Range("I6").Value = [double]"4.7514536831624372E-6"
$ws.Range("J6").Value = [double]"2.146518400025462E-3"
$ws.Range("K6").Value = [double]"0.90595022661523561"
$ws.Range("L6").Value = [double]"-0.55476895706741303"
$ws.Range("M6").Value = [double]"-0.60928128943968973"
$ws.Range("A7").Value = [double]"4.3949326748739852E-2"
$ws.Range("B7").Value = [double]"0.52029836992906986"
$ws.Range("C7").Value = [double]"3.5127929606025433E-2"
$ws.Range("D7").Value = [double]"167031315.02645785"
$ws.Range("E7").Value = [double]"38.921159788497206"
$ws.Range("F7").Value = [double]"54.514373093130857"
$ws.Range("G7").Value = [double]"3.3220014007090256E-12"
$ws.Range("H7").Value = [double]"0.24370926095811904"
$ws.Range("I7").Value = [double]"0.41778742603721808"
$ws.Range("J7").Value = [double]"4.9692977436049645E-7"
$ws.Range("K7").Value = [double]"0.99702926510373391"
$ws.Range("L7").Value = [double]"0.94411591506115933"
$ws.Range("M7").Value = [double]"0.79693513852082543"
$ws.Range("A8").Value = [double]"1.000000000223847E-4"
$ws.Range("B8").Value = [double]"9.975837554356189E-3"
$ws.Range("C8").Value = [double]"4.5393023297652463E-2"
$ws.Range("D8").Value = [double]"150000104.14826077"
$ws.Range("E8").Value = [double]"40.694482850043144"
$ws.Range("F8").Value = [double]"54.0868539502628"
$ws.Range("G8").Value = [double]"2.7154321805668484E-14"
$ws.Range("H8").Value = [double]"0.59999539377279199"
$ws.Range("I8").Value = [double]"1.8921872581286608E-5"
$ws.Range("J8").Value = [double]"2.1842029496302222E-3"
$ws.Range("K8").Value = [double]"0.84843861247920316"
$ws.Range("L8").Value = [double]"-0.49777429643210525"
$ws.Range("M8").Value = [double]"-0.5485732345049934"
$ws.Range("A9").Value = [double]"1.0000000002221964E-4"
$ws.Range("B9").Value = [double]"9.99891074427855E-3"
$ws.Range("C9").Value = [double]"4.5451571721433881E-2"
$ws.Range("D9").Value = [double]"150000016.99057075"
$ws.Range("E9").Value = [double]"43.365238605671117"
$ws.Range("F9").Value = [double]"56.082561606948843"
$ws.Range("G9").Value = [double]"2.2500527573449023E-14"
$ws.Range("H9").Value = [double]"0.59999975041590137"
$ws.Range("I9").Value = [double]"1.1673589495303797E-6"
$ws.Range("J9").Value = [double]"2.2349114779343406E-3"
$ws.Range("K9").Value = [double]"0.26577405463025205"
$ws.Range("L9").Value = [double]"-0.35818079583001161"
$ws.Range("M9").Value = [double]"-0.44984185955055822"
$ws.Range("A10").Value = [double]"7.5963745155290113E-3"
$ws.Range("B10").Value = [double]"1.2206981570078006"
$ws.Range("C10").Value = [double]"2.8393264414133448E-2"
$ws.Range("D10").Value = [double]"100000821.91427381"
$ws.Range("E10").Value = [double]"40.308858307296759"
$ws.Range("F10").Value = [double]"45.861540984734944"
$ws.Range("G10").Value = [double]"1.4230610827107152E-3"
$ws.Range("H10").Value = [double]"9.364078110130427E-3"
$ws.Range("I10").Value = [double]"0.89776261878227415"
$ws.Range("J10").Value = [double]"9.990725774741227E-2"
$ws.Range("K10").Value = [double]"0.99116142296851706"
$ws.Range("L10").Value = [double]"0.99277120216525649"
$ws.Range("M10").Value = [double]"0.96101635938080687"
$ws.Range("A11").Value = [double]"3.5288910092637994E-2"
$ws.Range("B11").Value = [double]"0.71282465918158855"
$ws.Range("C11").Value = [double]"3.0373481170051959E-2"
$ws.Range("D11").Value = [double]"158506227.79992354"
$ws.Range("E11").Value = [double]"30.000281635456471"
$ws.Range("F11").Value = [double]"59.448015374730666"
$ws.Range("G11").Value = [double]"2.0311374210082391"
$ws.Range("H11").Value = [double]"0.70163278652370376"
$ws.Range("I11").Value = [double]"0.6761632324664526"
$ws.Range("J11").Value = [double]"1.7736656090081997E-5"
$ws.Range("K11").Value = [double]"0.94492216471923152"
$ws.Range("L11").Value = [double]"0.95345550943150337"
$ws.Range("M11").Value = [double]"0.85268035963861666"
$ws.Range("A12").Value = [double]"1.0000000002228757E-4"
$ws.Range("B12").Value = [double]"9.9911537807638624E-3"
$ws.Range("C12").Value = [double]"4.5428461085543645E-2"
$ws.Range("D12").Value = [double]"150000052.37963778"
$ws.Range("E12").Value = [double]"39.189724282162636"
$ws.Range("F12").Value = [double]"54.945978383292939"
$ws.Range("G12").Value = [double]"2.404504331169536E-14"
$ws.Range("H12").Value = [double]"0.59999753406888989"
$ws.Range("I12").Value = [double]"6.5648983351732598E-6"
$ws.Range("J12").Value = [double]"2.1994711996414407E-3"
$ws.Range("K12").Value = [double]"0.28557380543756561"
$ws.Range("L12").Value = [double]"-0.5620710242423459"
$ws.Range("M12").Value = [double]"-0.67242711449880543"
$ws.Range("A13").Value = [double]"1.0000000002225485E-4"
$ws.Range("B13").Value = [double]"9.991967765815948E-3"
$ws.Range("C13").Value = [double]"4.5436017122414338E-2"
$ws.Range("D13").Value = [double]"150000057.54210895"
$ws.Range("E13").Value = [double]"42.31121219943018"
$ws.Range("F13").Value = [double]"54.851074934539319"
$ws.Range("G13").Value = [double]"2.3252119215215031E-14"
$ws.Range("H13").Value = [double]"0.59997235181215836"
$ws.Range("I13").Value = [double]"6.1517232934883539E-6"
$ws.Range("J13").Value = [double]"2.2213845843792749E-3"
$ws.Range("K13").Value = [double]"0.32840793534052026"
$ws.Range("L13").Value = [double]"-0.54309227761579648"
$ws.Range("M13").Value = [double]"-0.60270233357793934"
$ws.Range("A14").Value = [double]"1.0000000002226381E-4"
$ws.Range("B14").Value = [double]"9.9895163934579298E-3"
$ws.Range("C14").Value = [double]"4.5430826392854173E-2"
$ws.Range("D14").Value = [double]"150000041.51164338"
$ws.Range("E14").Value = [double]"41.326522744433163"
$ws.Range("F14").Value = [double]"54.044763907877787"
$ws.Range("G14").Value = [double]"2.345945196305314E-14"
$ws.Range("H14").Value = [double]"0.59998340305728659"
$ws.Range("I14").Value = [double]"7.6018322802250346E-6"
$ws.Range("J14").Value = [double]"2.210087869255555E-3"
$ws.Range("K14").Value = [double]"0.30552119362089014"
$ws.Range("L14").Value = [double]"-0.58009754217083165"
$ws.Range("M14").Value = [double]"-0.61782349358467181"
$ws.Range("A15").Value = [double]"1.0000000002226046E-4"
$ws.Range("B15").Value = [double]"9.9930054925361465E-3"
$ws.Range("C15").Value = [double]"4.5434409547244231E-2"
$ws.Range("D15").Value = [double]"150000044.5235678"
$ws.Range("E15").Value = [double]"41.912374303423633"
$ws.Range("F15").Value = [double]"55.190921131135397"
$ws.Range("G15").Value = [double]"2.3381066737607774E-14"
$ws.Range("H15").Value = [double]"0.59999822620323862"
$ws.Range("I15").Value = [double]"5.5066705438810774E-6"
$ws.Range("J15").Value = [double]"2.2187056496047229E-3"
$ws.Range("K15").Value = [double]"0.29411523227446279"
$ws.Range("L15").Value = [double]"-0.6419429321962038"
$ws.Range("M15").Value = [double]"-0.69907497269129593"
$ws.Range("A16").Value = [double]"1.0000000002222357E-4"
$ws.Range("B16").Value = [double]"1.0013624284674887E-2"
$ws.Range("C16").Value = [double]"4.5451775021870866E-2"
$ws.Range("D16").Value = [double]"150000075.63524768"
$ws.Range("E16").Value = [double]"45.938614337129245"
$ws.Range("F16").Value = [double]"57.994752503347101"
$ws.Range("G16").Value = [double]"2.2579824910258041E-14"
$ws.Range("H16").Value = [double]"0.59999966324973553"
$ws.Range("I16").Value = [double]"1.5548494452954616E-6"
$ws.Range("J16").Value = [double]"2.2507764940887043E-3"
$ws.Range("K16").Value = [double]"0.58158663110414355"
$ws.Range("L16").Value = [double]"-0.36120289674096839"
$ws.Range("M16").Value = [double]"-0.47511415946992241"
